$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the greeting text for the "R10" rule row from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Make E8 the active/selected cell, as it was the last cell edited
$ws.Activate()
$ws.Range("E8").Select()
